# Generate Report for Handback
# Adds a second handback row (new GUID file) to each sheet of the
# handback-status workbook, alongside refreshed timestamps / hashes for
# the existing (first) row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Constants
# ---------------------------------------------------------------------
$guid1 = "2d7a01be-cf30-4d1c-8039-d0c3be8f9467"   # existing row, new id
$guid2 = "56aee231-c2b1-473d-91ce-5a9c99f70f35"   # brand new row

$hash1 = "32c6422d0adc2da7f91cf13e0fbc88425583cb80"
$hash2 = "a38992c6af6ce8dba28ffe949505cd923eac2a9c"

$genDate        = "2016-08-29 09:05:40"
$zhHandoffDate  = "2016-08-29 09:05:35"
$zhHandbackDate = "2016-08-29 09:05:52"
$deHandbackDate = "2016-08-29 09:05:59"

$handedBack = "Handed back: in sync with en-US"

$md1 = "$guid1.md"
$md2 = "$guid2.md"
$e2eMd1 = "e2e\$guid1.md"
$e2eMd2 = "e2e\$guid2.md"

$xlf1zh = "$guid1.$hash1.zh-cn.xlf"
$xlf1de = "$guid1.$hash1.de-de.xlf"
$xlf2zh = "$guid2.$hash2.zh-cn.xlf"
$xlf2de = "$guid2.$hash2.de-de.xlf"

$urlBase1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/03703d013bfe43f59abc7adba445dd4c9650c500/e2e/$md1"
$urlBase2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/03703d013bfe43f59abc7adba445dd4c9650c500/e2e/$md2"
$urlZhCn1 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4738ca8993cd8952f6b8852a84c96b5d147f8cc1/e2e/$md1"
$urlZhCn2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4738ca8993cd8952f6b8852a84c96b5d147f8cc1/e2e/$md2"
$urlDeDe1 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/56b4e99d39abcba347ddcdb52686774d46889055/e2e/$md1"
$urlDeDe2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/56b4e99d39abcba347ddcdb52686774d46889055/e2e/$md2"

$hlinkColor = 15570276   # RGB(0x64,0x95,0xED) == the existing "HyperLink" cornflower blue

function Style-AsHyperlink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hlinkColor
}

# =======================================================================
# Sheet "Overview"
# =======================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

# --- row 2 : refresh the existing handback entry ---
$wsOverview.Range("A2").Value = $md1
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Range("B2").Value = $e2eMd1
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $urlBase1, "", "", $e2eMd1)
Style-AsHyperlink $wsOverview.Range("B2")
$wsOverview.Range("C2").Value = ".md"
$wsOverview.Range("E2").Value = $handedBack
$wsOverview.Range("F2").Value = $handedBack
$wsOverview.Range("G2").Value = $genDate

# --- row 3 : brand new handback entry ---
$wsOverview.Range("A3").Value = $md2
$wsOverview.Range("B3").Value = $e2eMd2
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $urlBase2, "", "", $e2eMd2)
Style-AsHyperlink $wsOverview.Range("B3")
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = $handedBack
$wsOverview.Range("F3").Value = $handedBack
$wsOverview.Range("G3").Value = $genDate
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# =======================================================================
# Sheet "zh-cn"
# =======================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

# --- row 2 : refresh the existing handback entry ---
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("A2").Value = $md1
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $urlBase1, "", "", $md1)
Style-AsHyperlink $wsZh.Range("A2")

$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = $handedBack
$wsZh.Range("D2").Value = "e2e"
$wsZh.Range("E2").Value = "ht"
$wsZh.Range("F2").Value = "False"
$wsZh.Range("G2").Value = $xlf1zh
$wsZh.Range("H2").Value = $zhHandoffDate
$wsZh.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Range("I2").Hyperlinks.Delete()
$wsZh.Range("I2").Value = $md1
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlZhCn1, "", "", $md1)
Style-AsHyperlink $wsZh.Range("I2")

$wsZh.Range("J2").Value = $xlf1zh
$wsZh.Range("K2").Value = $zhHandbackDate
$wsZh.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L2").Value = ""
$wsZh.Range("M2").Value = "True"
$wsZh.Range("N2").Value = ""
$wsZh.Range("O2").Value = "False"
$wsZh.Range("P2").Value = ""

# --- row 3 : brand new handback entry ---
$wsZh.Range("A3").Value = $md2
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $urlBase2, "", "", $md2)
Style-AsHyperlink $wsZh.Range("A3")

$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $handedBack
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = $xlf2zh
$wsZh.Range("H3").Value = $zhHandoffDate
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Range("I3").Value = $md2
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlZhCn2, "", "", $md2)
Style-AsHyperlink $wsZh.Range("I3")

$wsZh.Range("J3").Value = $xlf2zh
$wsZh.Range("K3").Value = $zhHandbackDate
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# =======================================================================
# Sheet "de-de"
# =======================================================================
$wsDe = $wb.Worksheets.Item("de-de")

# --- row 2 : refresh the existing handback entry ---
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("A2").Value = $md1
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $urlBase1, "", "", $md1)
Style-AsHyperlink $wsDe.Range("A2")

$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = $handedBack
$wsDe.Range("D2").Value = "e2e"
$wsDe.Range("E2").Value = "ht"
$wsDe.Range("F2").Value = "False"
$wsDe.Range("G2").Value = $xlf1de
$wsDe.Range("H2").Value = $genDate
$wsDe.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Range("I2").Hyperlinks.Delete()
$wsDe.Range("I2").Value = $md1
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlDeDe1, "", "", $md1)
Style-AsHyperlink $wsDe.Range("I2")

$wsDe.Range("J2").Value = $xlf1de
$wsDe.Range("K2").Value = $deHandbackDate
$wsDe.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L2").Value = ""
$wsDe.Range("M2").Value = "True"
$wsDe.Range("N2").Value = ""
$wsDe.Range("O2").Value = "False"
$wsDe.Range("P2").Value = ""

# --- row 3 : brand new handback entry ---
$wsDe.Range("A3").Value = $md2
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $urlBase2, "", "", $md2)
Style-AsHyperlink $wsDe.Range("A3")

$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $handedBack
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = $xlf2de
$wsDe.Range("H3").Value = $genDate
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Range("I3").Value = $md2
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlDeDe2, "", "", $md2)
Style-AsHyperlink $wsDe.Range("I3")

$wsDe.Range("J3").Value = $xlf2de
$wsDe.Range("K3").Value = $deHandbackDate
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))

Write-Output "done"
